# "Add cantrals by cantons"
#
# Sheet1 currently has a 2-row header (row 1 holds unit labels in E/G/I/J/K,
# row 2 holds "Hiver"/"Eté"/"Année" sub-labels) above 3 data rows (rows 3-5:
# Kembs, Birsfelden, Kembs-Centrale de dotation 1).
#
# The target layout instead uses a single header row with explicit column
# names (idx, idx2, Name, Date Start, Date End, (m3/s), (MW1), (MW2),
# (GWh) Winter, (GWh) Summer, (GWh) Year) directly above the 3 data rows
# (which then become rows 2-4).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the old two header rows — the three data rows (Kembs, Birsfelden,
# Kembs-Centrale de dotation 1) shift up and land on rows 1-3.
$ws.Rows.Item(1).Delete()
$ws.Rows.Item(1).Delete()

# Open up a fresh row above the data again to hold the new single-row
# header.
$ws.Rows.Item(1).Insert()

# Plain (unstyled) header cells.
$ws.Cells.Item(1,1).Value = "idx"
$ws.Cells.Item(1,2).Value = "idx2"
$ws.Cells.Item(1,3).Value = "Name"
$ws.Cells.Item(1,4).Value = "Date Start"
$ws.Cells.Item(1,5).Value = "Date End"

# Remaining header cells use the same Arial 9 font as the rest of the data
# table (matching the data-row name/unit cells). Columns F..K, in order.
$unitHeaders = "(m3/s)", "(MW1)", "(MW2)", "(GWh) Winter", "(GWh) Summer", "(GWh) Year"
$col = 6
foreach ($label in $unitHeaders) {
    $cell = $ws.Cells.Item(1, $col)
    $cell.Font.Name = "Arial"
    $cell.Font.Size = 9
    $cell.Value = $label
    $col = $col + 1
}

# Match the workbook's resulting active selection (first data row).
$ws.Range("A2:K2").Select()
